$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,13
$arr[0,0] = 3.952283258933164
$arr[0,1] = 0.1296783303549205
$arr[0,2] = 0.4128363259067811
$arr[0,3] = 0.09055545280074817
$arr[0,4] = 0
$arr[0,5] = 0.002801644214518938
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0.01658575960190589
$arr[0,9] = 3.316374157980533
$arr[0,10] = 0.5501803914211365
$arr[0,11] = 0.8499319006861441
$arr[0,12] = 8.165284513734179
$arr[1,0] = 3.919779709163947
$arr[1,1] = 0.1241628124187883
$arr[1,2] = 0.4123556489398794
$arr[1,3] = 0.09086679630885008
$arr[1,4] = 0
$arr[1,5] = 0.002807480157546356
$arr[1,6] = 0
$arr[1,7] = 0
$arr[1,8] = 0.01595156025396527
$arr[1,9] = 3.279550891548979
$arr[1,10] = 0.5496944244020483
$arr[1,11] = 0.8451357118398874
$arr[1,12] = 8.092704647791493
$arr[2,0] = 3.90172865468341
$arr[2,1] = 0.1208596785845941
$arr[2,2] = 0.4122216901707532
$arr[2,3] = 0.0910780728974947
$arr[2,4] = 0
$arr[2,5] = 0.00281125188906766
$arr[2,6] = 0
$arr[2,7] = 0
$arr[2,8] = 0.01555603837699771
$arr[2,9] = 3.258633361046265
$arr[2,10] = 0.549610162436224
$arr[2,11] = 0.8425781711422715
$arr[2,12] = 8.048630425926319
$arr[3,0] = 3.894851773449716
$arr[3,1] = 0.1195344658118103
$arr[3,2] = 0.4122076578791649
$arr[3,3] = 0.09116923855043346
$arr[3,4] = 0
$arr[3,5] = 0.002812836447537163
$arr[3,6] = 0
$arr[3,7] = 0
$arr[3,8] = 0.01539330295693553
$arr[3,9] = 3.250534314779969
$arr[3,10] = 0.5496296987278839
$arr[3,11] = 0.8416333313432958
$arr[3,12] = 8.030792341145144
$arr[4,0] = 3.893738799751986
$arr[4,1] = 0.1193156700963272
$arr[4,2] = 0.4122077782057403
$arr[4,3] = 0.09118468302503135
$arr[4,4] = 0
$arr[4,5] = 0.002813102438798022
$arr[4,6] = 0
$arr[4,7] = 0
$arr[4,8] = 0.01536618630131237
$arr[4,9] = 3.249215130910812
$arr[4,10] = 0.5496361974991544
$arr[4,11] = 0.8414823223546435
$arr[4,12] = 8.027837719570897
$arr[5,0] = 3.901633971285776
$arr[5,1] = 0.1208417220931466
$arr[5,2] = 0.4122213366739231
$arr[5,3] = 0.091079281850309
$arr[5,4] = 0
$arr[5,5] = 0.002811273066102433
$arr[5,6] = 0
$arr[5,7] = 0
$arr[5,8] = 0.01555384999792153
$arr[5,9] = 3.258522414224785
$arr[5,10] = 0.5496102077351708
$arr[5,11] = 0.8425650344462667
$arr[5,12] = 8.048389360228668
$arr[6,0] = 3.940680200131681
$arr[6,1] = 0.1277591831280063
$arr[6,2] = 0.4126371507262263
$arr[6,3] = 0.09065863715477374
$arr[6,4] = 0
$arr[6,5] = 0.00280361743547976
$arr[6,6] = 0
$arr[6,7] = 0
$arr[6,8] = 0.01636834714103941
$arr[6,9] = 3.303325986838274
$arr[6,10] = 0.5499684050434794
$arr[6,11] = 0.8481977653880435
$arr[6,12] = 8.140156793295091
$arr[7,0] = 4.032397241760975
$arr[7,1] = 0.1419931556691552
$arr[7,2] = 0.4147307766628501
$arr[7,3] = 0.08999281629069333
$arr[7,4] = 0
$arr[7,5] = 0.00279009248910179
$arr[7,6] = 0
$arr[7,7] = 0
$arr[7,8] = 0.01791786871955736
$arr[7,9] = 3.404642897330149
$arr[7,10] = 0.552369154997578
$arr[7,11] = 0.8623196122187338
$arr[7,12] = 8.324043949240632
$arr[8,0] = 4.109061144341638
$arr[8,1] = 0.1528695793555244
$arr[8,2] = 0.4170480962009577
$arr[8,3] = 0.08959994715350561
$arr[8,4] = 0
$arr[8,5] = 0.002781052219978708
$arr[8,6] = 0
$arr[8,7] = 0
$arr[8,8] = 0.01902849221183089
$arr[8,9] = 3.487341488136622
$arr[8,10] = 0.5551684916009947
$arr[8,11] = 0.8745765537953787
$arr[8,12] = 8.461616241506078
$arr[9,0] = 4.145963521021713
$arr[9,1] = 0.157910946748558
$arr[9,2] = 0.4182715475174632
$arr[9,3] = 0.08944199746583692
$arr[9,4] = 0
$arr[9,5] = 0.002777132005531268
$arr[9,6] = 0
$arr[9,7] = 0
$arr[9,8] = 0.01952798437619663
$arr[9,9] = 3.526770639874485
$arr[9,10] = 0.5566669772687476
$arr[9,11] = 0.8805627496737216
$arr[9,12] = 8.524754685357948
$arr[10,0] = 4.160229737542863
$arr[10,1] = 0.1598336236881153
$arr[10,2] = 0.4187591751703934
$arr[10,3] = 0.08938516147386721
$arr[10,4] = 0
$arr[10,5] = 0.002775674997118719
$arr[10,6] = 0
$arr[10,7] = 0
$arr[10,8] = 0.01971632512633548
$arr[10,9] = 3.541962475806315
$arr[10,10] = 0.557266776048067
$arr[10,11] = 0.8828886830557749
$arr[10,12] = 8.548744594714321
$arr[11,0] = 4.157144253246145
$arr[11,1] = 0.1594189325026321
$arr[11,2] = 0.4186530737560332
$arr[11,3] = 0.08939726992100816
$arr[11,4] = 0
$arr[11,5] = 0.002775987569637373
$arr[11,6] = 0
$arr[11,7] = 0
$arr[11,8] = 0.01967579812763631
$arr[11,9] = 3.538679028023466
$arr[11,10] = 0.5571361597566806
$arr[11,11] = 0.8823851224518222
$arr[11,12] = 8.543574334514176
$arr[12,0] = 4.14713135427445
$arr[12,1] = 0.1580688526526046
$arr[12,2] = 0.4183111773190404
$arr[12,3] = 0.08943726194728718
$arr[12,4] = 0
$arr[12,5] = 0.002777011586600574
$arr[12,6] = 0
$arr[12,7] = 0
$arr[12,8] = 0.01954349533059485
$arr[12,9] = 3.52801524888082
$arr[12,10] = 0.5567156746720201
$arr[12,11] = 0.8807529209312293
$arr[12,12] = 8.526726726763286
$arr[13,0] = 4.141036218963279
$arr[13,1] = 0.1572436682618843
$arr[13,2] = 0.4181049244942301
$arr[13,3] = 0.08946214549296272
$arr[13,4] = 0
$arr[13,5] = 0.002777642402179791
$arr[13,6] = 0
$arr[13,7] = 0
$arr[13,8] = 0.01946235162899868
$arr[13,9] = 3.521517371281277
$arr[13,10] = 0.5564623287802419
$arr[13,11] = 0.8797608478784866
$arr[13,12] = 8.516417620055108
$arr[14,0] = 4.106690331470247
$arr[14,1] = 0.1525420115250995
$arr[14,2] = 0.4169715463972068
$arr[14,3] = 0.0896106865633044
$arr[14,4] = 0
$arr[14,5] = 0.00278131227065477
$arr[14,6] = 0
$arr[14,7] = 0
$arr[14,8] = 0.01899573539953181
$arr[14,9] = 3.484801165804924
$arr[14,10] = 0.5550750899796242
$arr[14,11] = 0.8741936068820024
$arr[14,12] = 8.457501252859515
$arr[15,0] = 4.086139917788842
$arr[15,1] = 0.1496817885283974
$arr[15,2] = 0.4163196069861357
$arr[15,3] = 0.08970712358716604
$arr[15,4] = 0
$arr[15,5] = 0.002783612746464328
$arr[15,6] = 0
$arr[15,7] = 0
$arr[15,8] = 0.01870802561606411
$arr[15,9] = 3.462740831486542
$arr[15,10] = 0.554281699267122
$arr[15,11] = 0.8708834583621368
$arr[15,12] = 8.42150095855385
$arr[16,0] = 4.074510676034208
$arr[16,1] = 0.148045467042067
$arr[16,2] = 0.4159605636821198
$arr[16,3] = 0.08976454704810877
$arr[16,4] = 0
$arr[16,5] = 0.002784954023879151
$arr[16,6] = 0
$arr[16,7] = 0
$arr[16,8] = 0.01854200152723351
$arr[16,9] = 3.450222551553338
$arr[16,10] = 0.5538465453113588
$arr[16,11] = 0.8690181777167254
$arr[16,12] = 8.400846838655468
$arr[17,0] = 4.070605971357168
$arr[17,1] = 0.1474929439109474
$arr[17,2] = 0.415841734993748
$arr[17,3] = 0.08978432581081996
$arr[17,4] = 0
$arr[17,5] = 0.002785411271413802
$arr[17,6] = 0
$arr[17,7] = 0
$arr[17,8] = 0.0184856951316803
$arr[17,9] = 3.446013294650697
$arr[17,10] = 0.5537028481998618
$arr[17,11] = 0.8683932583471119
$arr[17,12] = 8.393862666935945
$arr[18,0] = 4.088307792327953
$arr[18,1] = 0.1499853519304395
$arr[18,2] = 0.4163873579202431
$arr[18,3] = 0.08969665539662763
$arr[18,4] = 0
$arr[18,5] = 0.002783365984276483
$arr[18,6] = 0
$arr[18,7] = 0
$arr[18,8] = 0.01873870867305882
$arr[18,9] = 3.465071567337873
$arr[18,10] = 0.5543639647252121
$arr[18,11] = 0.8712318308663498
$arr[18,12] = 8.425327833144252
$arr[19,0] = 4.150064453934704
$arr[19,1] = 0.1584650329092483
$arr[19,2] = 0.4184109402702632
$arr[19,3] = 0.0894254346302148
$arr[19,4] = 0
$arr[19,5] = 0.002776710063091518
$arr[19,6] = 0
$arr[19,7] = 0
$arr[19,8] = 0.0195823775966204
$arr[19,9] = 3.531140373451365
$arr[19,10] = 0.5568383032835555
$arr[19,11] = 0.8812307337993843
$arr[19,12] = 8.531673079714437
$arr[20,0] = 4.192128774742343
$arr[20,1] = 0.1640864602659065
$arr[20,2] = 0.419875297135718
$arr[20,3] = 0.08926551862386845
$arr[20,4] = 0
$arr[20,5] = 0.002772520205481997
$arr[20,6] = 0
$arr[20,7] = 0
$arr[20,8] = 0.02012907428054334
$arr[20,9] = 3.575841072545927
$arr[20,10] = 0.5586440132258446
$arr[20,11] = 0.8881100576827947
$arr[20,12] = 8.601646893559519
$arr[21,0] = 4.169522270709479
$arr[21,1] = 0.16107887499345
$arr[21,2] = 0.4190807680149078
$arr[21,3] = 0.08934928535418329
$arr[21,4] = 0
$arr[21,5] = 0.002774741806185552
$arr[21,6] = 0
$arr[21,7] = 0
$arr[21,8] = 0.01983771485456032
$arr[21,9] = 3.551844059425832
$arr[21,10] = 0.5576630180190847
$arr[21,11] = 0.8844068906299682
$arr[21,12] = 8.564257190611215
$arr[22,0] = 4.087327118606197
$arr[22,1] = 0.1498480858269318
$arr[22,2] = 0.4163566786166655
$arr[22,3] = 0.08970138189589072
$arr[22,4] = 0
$arr[22,5] = 0.00278347748705272
$arr[22,6] = 0
$arr[22,7] = 0
$arr[22,8] = 0.01872483878267417
$arr[22,9] = 3.464017329331739
$arr[22,10] = 0.5543267071722795
$arr[22,11] = 0.8710742140238636
$arr[22,12] = 8.423597569486674
$arr[23,0] = 4.005959150354272
$arr[23,1] = 0.138069758723816
$arr[23,2] = 0.4140275705303083
$arr[23,3] = 0.09015597843160794
$arr[23,4] = 0
$arr[23,5] = 0.002793593157406062
$arr[23,6] = 0
$arr[23,7] = 0
$arr[23,8] = 0.01750365813197519
$arr[23,9] = 3.375787488757709
$arr[23,10] = 0.5515378632368311
$arr[23,11] = 0.8581693247497597
$arr[23,12] = 8.273869122516402
$ws.Range("B2:N25").Value = $arr
